$d = $word.ActiveDocument

# 1. Drop the trailing period from the "pickup" line.
$d.Content.Find.Execute(
    "pick an item up when one is present in room.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "pick an item up when one is present in room", 2) | Out-Null

# 2. Insert a new "catch = ..." line right after the _GoBack bookmark
#    (between it and the "Quit = exit the game" line), matching the
#    line-break + text shape of its sibling runs.
$newText = "catch = catch when something in a room is catchable"

$rng = $d.Content
$rng.Find.Execute("Quit = exit the game", $true, $false, $false, $false,
                   $false, $true, 1, $false, "", 0) | Out-Null
$quitStart = $rng.Start

# Insert "<text><line break>" right before the "Quit" run's own text; the
# line break that used to sit immediately in front of "Quit" ends up in
# front of "catch" instead (reused), and the break we append becomes the
# new lead-in for "Quit" - so nothing behind the bookmark needs to move.
$insertPoint = $d.Range($quitStart, $quitStart)
$insertPoint.InsertBefore($newText + [char]11)

# Force the newly inserted "<break><catch text>" span to become its own
# run (distinct from both the preceding "pickup" run and the following
# "Quit" run) by round-tripping a character formatting property over
# exactly that span.
$breakPos = $quitStart - 1
$catchRngEnd = $breakPos + 1 + $newText.Length
$catchRng = $d.Range($breakPos, $catchRngEnd)
$catchRng.Font.Bold = 1
$catchRng.Font.Bold = 0
